$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo in the "login using ..." string (row 10, column B)
$ws.Range("B10").Value = "login using email / username, and password"

# Add the new "status" column (column A) header
$ws.Range("A1").Value = "status"

# Mark existing implemented endpoints as "ok" in the new status column
$ws.Range("A2").Value = "ok"
$ws.Range("A4").Value = "ok"
$ws.Range("A5").Value = "ok"
$ws.Range("A6").Value = "ok"

# "get user info" (row 6) API: widen permission to any logged user
$ws.Range("E6").Value = "any logged user"

# Move the active selection to A7
$null = $ws.Range("A7").Select()
